# PlayerDna.xlsx - "a new dna plan"
# Column I ("~Calculater") moves from a calculated POWER(2,A4) column to a
# hand-authored "归类" (classification) column of text descriptions, and one
# bad MutexId value is corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a bad MutexId value (row for Id=16 should mutex with 17 & 18, not 17 twice) ---
$ws.Range("E19").Value = "17,18"

# --- Column I header changes from "计算系数" (calc coefficient) to "归类" (category) ---
$ws.Range("I1").Value = "归类"
$ws.Range("I2").Value = "string"

# --- Replace the POWER(2,A4) calculated column with static classification text ---
$ws.Range("I4").Value = "冒险概率+ 交易概率+"
$ws.Range("I5").Value = "战斗难度- 角色概率-"
$ws.Range("I6").Value = "陷阱难度- 采集概率-"
$ws.Range("I7").Value = "鉴定难度- 冒险概率-"
$ws.Range("I8").Value = "隐藏的对话选项"
$ws.Range("I9").Value = "战斗难度- 交易概率-"
$ws.Range("I10").Value = "陷阱难度+ 宝物概率+"
$ws.Range("I11").Value = "战斗概率+ 冒险概率+"
$ws.Range("I12").Value = "隐藏的对话选项-偷袭"
$ws.Range("I13").Value = "隐藏的对话选项-单挑"
$ws.Range("I14").Value = "游戏概率+ 采集概率+"
$ws.Range("I15").Value = "增益难度- 战斗难度+"
$ws.Range("I16").Value = "战斗概率- 冒险难度+"
$ws.Range("I17").Value = "游戏难度+ 增益概率+"
$ws.Range("I18").Value = "增益难度- 游戏概率+"
$ws.Range("I19").Value = "隐藏的对话选项-时间"
$ws.Range("I20").Value = "隐藏的对话选项-绕路"
$ws.Range("I21").Value = "冒险难度- 增益概率-"
$ws.Range("I22").Value = "宝物概率- 战斗难度-"
$ws.Range("I23").Value = "增益难度+ 卡牌概率+"
$ws.Range("I24").Value = "战斗概率+ 冒险概率+"
$ws.Range("I25").Value = "角色难度+ 游戏难度-"
$ws.Range("I26").Value = "boss难度++ 宝物概率+"
$ws.Range("I27").Value = "困难事件概率+++ 角色概率+"

# --- Format column I data as text (was numeric/General) and widen it to fit ---
$ws.Range("I4:I27").NumberFormat = "@"
$ws.Columns("I").ColumnWidth = 23.125

# --- Active cell moves to I10 ---
$ws.Range("I10").Select()
